$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest crypto
# quotes. Price strings that look like plain numbers are written with a
# leading apostrophe (the normal Excel "force text" entry) so they stay
# text cells - exactly like the source data - instead of being silently
# parsed into numbers (which would also strip meaningful trailing zeros,
# e.g. "1.00" -> 1 or "0.100" -> 0.1).

$ws.Range("D2").Value = "40.706.17"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").Value = "2.377.80"
$ws.Range("E3").Value = "  -3.63%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'310.59"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").Value = "'87.31"
$ws.Range("E6").Value = "  -5.38%  "
$ws.Range("D7").Value = "'0.529"
$ws.Range("E7").Value = "  -4.13%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.496"
$ws.Range("E9").Value = "  -3.13%  "
$ws.Range("D10").Value = "'0.0847"
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("D11").Value = "'30.72"
$ws.Range("E11").Value = "  -6.79%  "
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "2.742.83"
$ws.Range("E13").Value = "  -3.67%  "
$ws.Range("D14").Value = "'6.52"
$ws.Range("E14").Value = "  -5.17%  "
$ws.Range("D15").Value = "'15.11"
$ws.Range("E15").Value = "  -2.33%  "
$ws.Range("D16").Value = "2.390.52"
$ws.Range("E16").Value = "  -3.41%  "
$ws.Range("D17").Value = "'0.763"
$ws.Range("E17").Value = "  -3.85%  "
$ws.Range("D18").Value = "40.671.28"
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("D19").Value = "0.0₃0913"
$ws.Range("E19").Value = "  -3.46%  "
$ws.Range("D20").Value = "'6.15"
$ws.Range("E20").Value = "  -4.56%  "
$ws.Range("D21").Value = "'68.67"
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("D22").Value = "'10.77"
$ws.Range("E22").Value = "  -4.62%  "
$ws.Range("D23").Value = "'234.77"
$ws.Range("E23").Value = "  -2.27%  "
$ws.Range("E24").Value = "  -5.99%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "'1.81"
$ws.Range("E26").Value = "  -7.33%  "
$ws.Range("D27").Value = "'23.82"
$ws.Range("E27").Value = "  -4.12%  "
$ws.Range("D28").Value = "'2.21"
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("D29").Value = "'9.33"
$ws.Range("E29").Value = "  -3.82%  "
$ws.Range("D30").Value = "'34.01"
$ws.Range("E30").Value = "  -6.31%  "
$ws.Range("D31").Value = "'154.08"
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").Value = "'5.24"
$ws.Range("E33").Value = "  -4.26%  "
$ws.Range("E34").Value = "  -4.26%  "
$ws.Range("D35").Value = "'2.42"
$ws.Range("E35").Value = "  -5.54%  "
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("D37").Value = "'16.04"
$ws.Range("E37").Value = "  -6.89%  "

# Rows 38/39: Kaspa and LidoDAOToken swapped ranking positions.
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.78"
$ws.Range("E38").Value = "  -4.10%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.100"
$ws.Range("E39").Value = "  -3.72%  "

$ws.Range("E40").Value = "  -7.29%  "
$ws.Range("D41").Value = "'3.82"
$ws.Range("E41").Value = "  -4.73%  "
$ws.Range("D42").Value = "'2.37"
$ws.Range("E42").Value = "  -6.78%  "
$ws.Range("D43").Value = "1.969.60"
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("D44").Value = "'0.0270"
$ws.Range("E44").Value = "  -4.78%  "
$ws.Range("D45").Value = "'17.85"
$ws.Range("E45").Value = "  -5.64%  "
$ws.Range("D46").Value = "'9.51"
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").Value = "'2.71"
$ws.Range("E47").Value = "  -8.41%  "
$ws.Range("D48").Value = "2.606.11"
$ws.Range("E48").Value = "  -3.67%  "
$ws.Range("D49").Value = "'93.24"
$ws.Range("E49").Value = "  -4.46%  "
$ws.Range("D50").Value = "'72.43"
$ws.Range("E50").Value = "  -3.94%  "
$ws.Range("D51").Value = "'50.49"
$ws.Range("E51").Value = "  -3.56%  "
